$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number issue bump, week-of dates) ---
$ws.Range("A8").Value = "Volume 33   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/19/2026  Through  1/25/2026"

# --- Row 16 (Robbery) ---
$ws.Range("F16").Value = 2
$ws.Range("I16").Value = 2
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = -71.428571428571

# --- Row 17 (Fel. Assault) ---
$ws.Range("F17").Value = 2

# --- Row 21 (Gr. Larceny) ---
$ws.Range("I21").Value = 2
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -33.333333333333
$ws.Range("N21").Value = -88.888888888888

# --- Row 24 (Housing) : cells flip from blank placeholder text to real numbers ---
$ws.Range("D24").Value = 1
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("E24").Value = -100
$ws.Range("E24").NumberFormat = "#,##0.0;""-""#,##0.0"

$ws.Range("G24").Value = 1
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("H24").Value = 0
$ws.Range("H24").NumberFormat = "#,##0.0;""-""#,##0.0"

$ws.Range("J24").Value = 1
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("K24").Value = -100
$ws.Range("K24").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Row 26 (Petit Larceny) ---
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 2
$ws.Range("K26").Value = 100
$ws.Range("M26").Value = 0

# --- Row 28 (Misd. Assault) ---
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("I28").Value = 1
$ws.Range("I28").NumberFormat = "#,##0"
